$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.351.70'
$ws.Range("E2").Value = '  +1.90%  '
$ws.Range("D3").Value = '1.882.21'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'0.697"
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").Value = "'246.52"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +5.41%  '
$ws.Range("D9").Value = "'0.356"
$ws.Range("E9").Value = '  +2.77%  '
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = "'13.58"
$ws.Range("E12").Value = '  +5.82%  '
$ws.Range("D13").Value = "'0.773"
$ws.Range("E13").Value = '  +8.52%  '
$ws.Range("D14").Value = '2.156.04'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '1.876.62'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '35.341.35'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").Value = "'73.42"
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").Value = '0.0₃0828'
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").Value = "'244.94"
$ws.Range("E20").Value = '  -0.66%  '
$ws.Range("D21").Value = "'12.80"
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("D22").Value = "'5.15"
$ws.Range("E22").Value = '  +4.92%  '
$ws.Range("D23").Value = "'2.62"
$ws.Range("E23").Value = '  +9.02%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  -3.95%  '
$ws.Range("D26").Value = "'164.66"
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = "'8.65"
$ws.Range("E27").Value = '  +3.29%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("E30").Value = '  +3.47%  '
$ws.Range("E31").Value = '  +1.33%  '
$ws.Range("B32").Value = 'WEMIXToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = "'1.88"
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = "'4.18"
$ws.Range("E33").Value = '  +1.08%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -11.83%  '
$ws.Range("D36").Value = "'0.854"
$ws.Range("E36").Value = '  +2.47%  '
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("E38").Value = '  +12.10%  '
$ws.Range("E39").Value = '  +4.49%  '
$ws.Range("D40").Value = "'17.23"
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").Value = "'97.19"
$ws.Range("E41").Value = '  -0.97%  '
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("D44").Value = '1.303.60'
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("D45").Value = "'0.0808"
$ws.Range("E45").Value = '  +5.55%  '
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").Value = "'12.01"
$ws.Range("E48").Value = '  -1.70%  '
$ws.Range("D50").Value = "'42.14"
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '2.060.50'
$ws.Range("E51").Value = '  +0.16%  '
